# Commit: "update ảnh cho 20 sản phẩm" (update images for 20 products)
#
# A batch of product image files under products2\ and products3\ no
# longer exist on disk, so the pic2 / pic3 cells that used to point at
# them are reset to "null" (the placeholder value already used
# elsewhere in this sheet for products that are missing that picture).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cellsToClear = @(
    "E67", "F67",
    "E69", "F69",
    "E75", "F75",
    "E76", "F76",
    "E77", "F77",
    "E78", "F78",
    "E79", "F79",
    "F80",
    "E81", "F81"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Value = "null"
}

# Reflect where the author was looking in the sheet after the edit.
$ws.Range("E81").Select()
